$wb = $excel.ActiveWorkbook

# --- Update values on L_map sheet (B2:B5: 1 -> 0, B6:B10: 2 -> 1) ---
$wsMap = $wb.Worksheets("L_map")
$wsMap.Range("B2:B5").Value = 0
$wsMap.Range("B6:B10").Value = 1

# --- Update view/selection on SS1 (becomes SS0) ---
$wsSS1 = $wb.Worksheets("SS1")
$null = $wsSS1.Range("B3").Select()
$excel.ActiveWindow.Zoom = 132

# --- Update view/selection on SS2 (becomes SS1) ---
$wsSS2 = $wb.Worksheets("SS2")
$null = $wsSS2.Range("H13").Select()

# --- Delete SS3 entirely ---
$excel.DisplayAlerts = $false
$null = $wb.Worksheets("SS3").Delete()
$excel.DisplayAlerts = $true

# --- Rename sheets: SS1 -> SS0, SS2 -> SS1 ---
$wsSS1.Name = "SS0"
$wsSS2.Name = "SS1"

# --- Update view/selection on L_map: drop tabSelected, move selection to C1 ---
$null = $wsMap.Range("C1").Select()

# --- Update view/selection on var: becomes the active/selected sheet with B2 selected ---
$wsVar = $wb.Worksheets("var")
$null = $wsVar.Range("B2").Select()

Write-Host "Done"
